# Jeannine's log.xlsx -- "Tweak to the ranking algo to better sort and
# distribute the events in a log." Adds a Monday (2016-09-19) and a
# Tuesday (2016-09-20) block of entries to the bottom of the Logs sheet,
# extends the Staff_Name named range / database sheet by one row, and
# updates the saved view state on both sheets.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$db   = $wb.Worksheets.Item("database")

# ---------------------------------------------------------------------
# 2. Extend the Staff_Name defined name to include the new row.
# ---------------------------------------------------------------------
$wb.Names.Item("Staff_Name").RefersTo = "=database!`$A`$2:`$A`$39"

# ---------------------------------------------------------------------
# 3. Logs sheet: append a MONDAY (9/19) section starting at row 199 and
#    a TUESDAY (9/20) section starting at row 212. Row 199/212 are
#    section-header rows (same banded format as the existing headers,
#    e.g. row 194) -- clone the format from row 194 via copy/paste
#    special so the style indices line up exactly.
# ---------------------------------------------------------------------
$logs.Range("A194:F194").Copy()
$logs.Range("A199:F199").PasteSpecial(-4122)  # xlPasteFormats
$logs.Range("A194:F194").Copy()
$logs.Range("A212:F212").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$logs.Range("B199").Value = "MONDAY"
$logs.Range("B212").Value = "TUESDAY"

function Set-LogRow($Row, $A, $B, $C, $D, $E, $F, $Height) {
    $logs.Range("A$Row").Value = $A
    $logs.Range("B$Row").Value = $B
    $logs.Range("C$Row").Value = $C
    $logs.Range("D$Row").Value = $D
    $logs.Range("E$Row").Value = $E
    $logs.Range("F$Row").Value = $F
    if ($Height) {
        $logs.Rows.Item($Row).RowHeight = $Height
    }
}

# ---- Monday 9/19/2016 (serial 42632) — Lassonde C wireless mic setup ----
Set-LogRow 200 "Setup Mic" 42632 "1715" "LAS" "C" `
    'Take cart with mixer, 2 wireless mics and 2 mic stands from Lassonde 1011 storeroom (across from Lassonde A). Go to Lassonde C classroom (class starts at 5:30 pm but be there early in case previous class ends early). ' `
    60

Set-LogRow 201 "Other" 42632 "1715" "LAS" "C" `
    'Log in as 5065*0 on touchscreen. (First level bar is your wireless handheld mic volume). Plug in mic cable from output of mixer to mic input on podium (XLR jack just above VHS machine in podium). Ramp up volume a bit on "Microphone 2" on touchscreen to medium volume to get level.' `
    75

Set-LogRow 202 "Other" 42632 "1715" "LAS" "C" `
    'Plug in power cord from cart on to power outlet on left side of podium (to left of document camera). Turn on mixer. Turn on wireless microphone receivers on cart (NOTE: DO NOT PRESS "SYNC" BUTTON" - POWER BUTTON IS FIRST BUTTON TO THE RIGHT ON RECEIVER). ' `
    75

Set-LogRow 203 "Other" 42632 "1715" "LAS" "C" `
    'Press "MUTE" button on wireless microphones to turn on mics. Adjust volume by adjusting volume on mixer (inputs 1 and 2). Also you can adjust volume on first volume control bar on touchscreen. (NOTE: VOLUME ON TOUCHSCREEN MUST BE RAMPED UP OR DOWN INITIALLY TO GET ANY VOLUME - First volume bar).' `
    90

Set-LogRow 204 "Other" 42632 "1715" "LAS" "C" `
    'Once volumes are set, place one mic stand with mic halfway up aisle on right and one mic stand with mic halfway up aisle on left. Demo volume controls to prof. and demo PC. Leave microphone bags with milk carton on cart in room. PLEASE FIND OUT END TIME OF CLASS FROM PROF. AND TELL MASI AS MICROPHONES ARE EXPENSIVE. TELL PROF. TO STAY WITH MICS UNTIL THEY ARE PICKED UP. TELL HIM TO CALL ext 55800   WHEN DONE (use phone in classroom).' `
    120

Set-LogRow 205 "Pickup Mic" 42632 "1850" "LAS" "C" `
    'Pick up 2 wireless mics on stands with cart. Move all equipment on cart - cart has 2 wireless mic receivers and mixer and mic cables. Pick up 2 mic stands - return all equipment to Lassonde 1011 storeroom (across the hall from Lassonde A). PLEASE PUT 2 WIRELESS MICS IN BAGS PROVIDED IN MILK CARTON ON CART. Very expensive mics - please go early and treat mics with care.' `
    90

Set-LogRow 206 "Other" 42632 "1850" "LAS" "C" `
    'Turn off wireless microphones by pressing "MUTE" button on mics.' `
    30

# NOTE: the shared-string table in the authored workbook picked up the
# F208 comment ("PLEASE BE ON TIME...") before the F207 comment ("Turn
# off wireless microphone receivers...") and before database!A39
# ("Null"). Poke that string in now (row 208 is filled in properly
# below) so new shared strings come out in the same order as the
# original edit.
$logs.Range("F208").Value = 'PLEASE BE ON TIME - Prof upset last week when no one came till 7:05 pm and other class was starting.'

Set-LogRow 207 "Operator" 42632 "1850" "LAS" "C" `
    'Turn off wireless microphone receivers by pressing "POWER" button and not "SYNC" button. ' `
    30

# ---------------------------------------------------------------------
# 1. database sheet: give A39 (previously blank) the value "Null" so the
#    Staff_Name list has one more non-blank row.
# ---------------------------------------------------------------------
$db.Range("A39").Value = "Null"

Set-LogRow 208 "Other" 42632 "1850" "LAS" "C" `
    'PLEASE BE ON TIME - Prof upset last week when no one came till 7:05 pm and other class was starting.' `
    30

# ---- Tuesday 9/20/2016 (serial 42633) ----
Set-LogRow 213 "Demo" 42633 "1600" "BC" "215" `
    "Make sure prof is happy. Go 10 minutes early to class." `
    0

Set-LogRow 214 "AV Shutdown" 42633 "1730" "BC" "230" `
    "Return keyboard and DVD remote to drawer." `
    0

# ---------------------------------------------------------------------
# 4. Update saved view state to match where the log was left scrolled.
# ---------------------------------------------------------------------
$logs.Activate()
[void]$logs.Range("E217").Select()
$excel.ActiveWindow.ScrollRow = 205

$db.Activate()
[void]$db.Range("A41").Select()

$logs.Activate()

Write-Host "Applied Jeannine's log update."
